# Fixed variables and query errors in Breed queries (TC01-TC30 series).
# This workbook (TC04_Canine_Filter_Breed-BassHnd) had an erroneous trailing
# "coalesce(co.cohort_description, '') AS `Cohort`" column appended to the
# "Cases" tab's Neo4j query (cell B2 on the startup sheet). That extra
# return column referenced a cohort (co) node the query never reliably
# matched for every case, so it is removed here - the query now ends
# cleanly at "Response to Treatment".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Correct the "CasesTab" query text in B2 (drop the trailing Cohort column) ---
$fixedCaseQuery = 'MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
WHERE demo.breed IN [''Basset Hound'']
MATCH (c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '''') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '''') AS `Study Code` ,
        coalesce(s.clinical_study_type, '''') AS  `Study Type`,
        coalesce(demo.breed, '''') AS Breed ,
        coalesce(diag.disease_term, '''') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '''') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '''') AS Age ,
        coalesce(demo.sex, '''') AS Sex ,
        coalesce(demo.neutered_indicator, '''') AS `Neutered Status`,
        coalesce(demo.weight, '''') AS `Weight (kg)`,
        coalesce(diag.best_response, '''') AS `Response to Treatment`
'

$ws.Range("B2").Value = $fixedCaseQuery

# --- Row heights shrink now that B2 wraps one fewer line (and the sheet's
#     font-metrics were re-measured on save) ---
$ws.Rows.Item(2).RowHeight = 259.2
$ws.Rows.Item(3).RowHeight = 230.4
$ws.Rows.Item(4).RowHeight = 244.8

# --- View state: zoom to 130%, and select the corrected B2 cell instead of
#     the old B4 selection ---
$ws.Activate()
$excel.ActiveWindow.Zoom = 130
$ws.Range("B2").Select()
